$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.672.10"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.175.97"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.84"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.99"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.174.67"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.23"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.111"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.400"
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("D13").Value = "3.729.01"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.03"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "58.734.50"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "3.177.45"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.98"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.10"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "357.14"
$ws.Range("E22").Value = "  +5.30%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.518"
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.59"
$ws.Range("E25").Value = "  +3.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").Value = "0.0₃0958"
$ws.Range("E27").Value = "  +4.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.59"
$ws.Range("E29").Value = "  +4.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.58"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.92"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.42"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.22"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.93"
$ws.Range("E35").Value = "  +6.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.91"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.28"
$ws.Range("E37").Value = "  +3.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.60"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.32"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("E40").Value = "  +14.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0678"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.07"
$ws.Range("E43").Value = "  +4.59%  "
$ws.Range("D44").Value = "3.219.13"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "37.05"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0274"
$ws.Range("E46").Value = "  +6.52%  "
$ws.Range("D47").Value = "2.351.43"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.02"
$ws.Range("E49").Value = "  +7.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.77"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.11"
$ws.Range("E51").Value = "  +1.84%  "
